# Applies the cryptos.xlsx data refresh described by the commit
# "Updated cryptos list on Thu Aug 22 14:53:19 UTC 2024 with GitHub Actions".
# Price (column D) and 1h volume (column E) values are refreshed for most rows,
# two coin/link columns are swapped (rows 31-32: Aptos <-> USDe), and a handful
# of price cells are textual numbers that must stay text (not be auto-converted
# to real numbers by Excel), so those are written via a brief Text-number-format
# round trip that restores each cell's original style afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

$ws.Range('D2').Value = '60.287.73'
$ws.Range('E2').Value = '  +1.50%  '
$ws.Range('D3').Value = '2.598.89'
$ws.Range('E3').Value = '  +0.61%  '
$ws.Range('E4').Value = '  -0.10%  '
Set-TextValue 'D5' '585.85'
$ws.Range('E5').Value = '  +5.98%  '
Set-TextValue 'D6' '142.72'
$ws.Range('E6').Value = '  +1.81%  '
$ws.Range('E7').Value = '  -0.06%  '
Set-TextValue 'D8' '0.598'
$ws.Range('E8').Value = '  +0.80%  '
$ws.Range('D9').Value = '2.609.10'
$ws.Range('E9').Value = '  +0.54%  '
Set-TextValue 'D10' '6.50'
$ws.Range('E10').Value = '  -3.38%  '
Set-TextValue 'D11' '0.105'
$ws.Range('E11').Value = '  +1.54%  '
$ws.Range('E12').Value = '  -2.78%  '
Set-TextValue 'D13' '0.369'
$ws.Range('E13').Value = '  +4.34%  '
$ws.Range('D14').Value = '3.065.50'
$ws.Range('E14').Value = '  +0.57%  '
Set-TextValue 'D15' '24.62'
$ws.Range('E15').Value = '  +6.62%  '
$ws.Range('D16').Value = '60.273.84'
$ws.Range('E16').Value = '  +1.43%  '
$ws.Range('E17').Value = '  +2.95%  '
$ws.Range('D18').Value = '2.608.68'
$ws.Range('E18').Value = '  +0.40%  '
$ws.Range('E19').Value = '  +9.81%  '
Set-TextValue 'D20' '4.65'
$ws.Range('E20').Value = '  +2.09%  '
Set-TextValue 'D21' '347.16'
$ws.Range('E21').Value = '  +2.33%  '
$ws.Range('E22').Value = '  +6.16%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('E24').Value = '  +9.32%  '
Set-TextValue 'D25' '63.15'
$ws.Range('E25').Value = '  -0.06%  '
Set-TextValue 'D26' '1.00'
$ws.Range('E26').Value = '  +0.39%  '
$ws.Range('E27').Value = '  +0.18%  '
Set-TextValue 'D28' '8.08'
$ws.Range('E28').Value = '  +8.16%  '
$ws.Range('D29').Value = '0.0₃0792'
$ws.Range('E29').Value = '  +2.50%  '
Set-TextValue 'D30' '1.87'
$ws.Range('E30').Value = '  +10.70%  '
$ws.Range('B31').Value = 'USDe'
$ws.Range('C31').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D31' '0.998'
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D32' '6.37'
$ws.Range('E32').Value = '  +3.79%  '
Set-TextValue 'D33' '163.16'
$ws.Range('E33').Value = '  +3.62%  '
$ws.Range('E34').Value = '  +1.79%  '
Set-TextValue 'D35' '4.27'
$ws.Range('E35').Value = '  +3.64%  '
Set-TextValue 'D36' '0.977'
$ws.Range('E36').Value = '  +8.76%  '
Set-TextValue 'D37' '1.23'
$ws.Range('E37').Value = '  +5.83%  '
$ws.Range('E38').Value = '  +10.35%  '
$ws.Range('E39').Value = '  +1.36%  '
$ws.Range('E40').Value = '  +6.46%  '
Set-TextValue 'D41' '307.75'
$ws.Range('E41').Value = '  +6.58%  '
Set-TextValue 'D42' '0.837'
$ws.Range('E42').Value = '  -0.38%  '
Set-TextValue 'D43' '135.65'
$ws.Range('E43').Value = '  -0.02%  '
Set-TextValue 'D44' '0.0993'
$ws.Range('E44').Value = '  +2.38%  '
$ws.Range('E45').Value = '  -0.02%  '
Set-TextValue 'D46' '19.78'
$ws.Range('E46').Value = '  +4.32%  '
Set-TextValue 'D47' '5.00'
$ws.Range('E47').Value = '  +10.45%  '
Set-TextValue 'D48' '0.603'
$ws.Range('E48').Value = '  +1.12%  '
Set-TextValue 'D49' '0.0549'
$ws.Range('E49').Value = '  +2.92%  '
Set-TextValue 'D50' '20.08'
$ws.Range('E50').Value = '  +7.75%  '
$ws.Range('E51').Value = '  +2.78%  '
